$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the last data row (old row 6) - data now spans rows 1-5
$ws.Rows.Item(6).Delete()

# 2) Replace the data values in rows 2-5 with the new dataset (1000-row export slice)
$ws.Cells.Item(2, 1).Value = 45093.50694444445
$ws.Cells.Item(2, 2).Value = 9.173
$ws.Cells.Item(2, 3).Value = 6.534
$ws.Cells.Item(2, 4).Value = 2.853
$ws.Cells.Item(2, 5).Value = 20.581
$ws.Cells.Item(2, 6).Value = 14.857
$ws.Cells.Item(2, 7).Value = 6.617
$ws.Cells.Item(2, 8).Value = 20.037
$ws.Cells.Item(2, 9).Value = 11.575
$ws.Cells.Item(2, 10).Value = 4.738
$ws.Cells.Item(2, 11).Value = 6.004
$ws.Cells.Item(2, 12).Value = 8.497
$ws.Cells.Item(2, 13).Value = 8.644
$ws.Cells.Item(2, 14).Value = 2.821
$ws.Cells.Item(2, 15).Value = 7.519
$ws.Cells.Item(2, 16).Value = 9.989000000000001
$ws.Cells.Item(2, 17).Value = 7.075
$ws.Cells.Item(2, 18).Value = 2.196
$ws.Cells.Item(2, 19).Value = 0.759
$ws.Cells.Item(2, 20).Value = 107.045
$ws.Cells.Item(2, 21).Value = 20.658
$ws.Cells.Item(2, 22).Value = 6.941
$ws.Cells.Item(2, 23).Value = 12.856
$ws.Cells.Item(2, 24).Value = 7.392
$ws.Cells.Item(2, 25).Value = 1.396
$ws.Cells.Item(2, 26).Value = 11.913
$ws.Cells.Item(2, 27).Value = 6.131
$ws.Cells.Item(2, 28).Value = 5.786
$ws.Cells.Item(2, 29).Value = 6.628
$ws.Cells.Item(2, 30).Value = 8.785
$ws.Cells.Item(2, 31).Value = 2.215
$ws.Cells.Item(2, 32).Value = 18.053
$ws.Cells.Item(2, 33).Value = 3.621
$ws.Cells.Item(2, 34).Value = 8.676
$ws.Cells.Item(3, 1).Value = 45093.51388888889
$ws.Cells.Item(3, 2).Value = 17.513
$ws.Cells.Item(3, 3).Value = 13.009
$ws.Cells.Item(3, 4).Value = 1.689
$ws.Cells.Item(3, 5).Value = 38.631
$ws.Cells.Item(3, 6).Value = 30.823
$ws.Cells.Item(3, 7).Value = 13.359
$ws.Cells.Item(3, 8).Value = 50.335
$ws.Cells.Item(3, 9).Value = 21.469
$ws.Cells.Item(3, 10).Value = 9.515000000000001
$ws.Cells.Item(3, 11).Value = 13.456
$ws.Cells.Item(3, 12).Value = 15.704
$ws.Cells.Item(3, 13).Value = 16.344
$ws.Cells.Item(3, 14).Value = 4.76
$ws.Cells.Item(3, 15).Value = 13.911
$ws.Cells.Item(3, 16).Value = 19.524
$ws.Cells.Item(3, 17).Value = 12.075
$ws.Cells.Item(3, 18).Value = 1.154
$ws.Cells.Item(3, 19).Value = 0.708
$ws.Cells.Item(3, 20).Value = 204.406
$ws.Cells.Item(3, 21).Value = 38.846
$ws.Cells.Item(3, 22).Value = 12.84
$ws.Cells.Item(3, 23).Value = 25.72
$ws.Cells.Item(3, 24).Value = 13.733
$ws.Cells.Item(3, 25).Value = 1.929
$ws.Cells.Item(3, 26).Value = 26.081
$ws.Cells.Item(3, 27).Value = 11.342
$ws.Cells.Item(3, 28).Value = 10.198
$ws.Cells.Item(3, 29).Value = 11.954
$ws.Cells.Item(3, 30).Value = 16.301
$ws.Cells.Item(3, 31).Value = 1.023
$ws.Cells.Item(3, 32).Value = 46.167
$ws.Cells.Item(3, 33).Value = 7.1
$ws.Cells.Item(3, 34).Value = 16.053
$ws.Cells.Item(4, 1).Value = 45093.52083333334
$ws.Cells.Item(4, 2).Value = 9.42
$ws.Cells.Item(4, 3).Value = 6.99
$ws.Cells.Item(4, 4).Value = 1.025
$ws.Cells.Item(4, 5).Value = 20.891
$ws.Cells.Item(4, 6).Value = 16.536
$ws.Cells.Item(4, 7).Value = 7.087
$ws.Cells.Item(4, 8).Value = 31.718
$ws.Cells.Item(4, 9).Value = 11.584
$ws.Cells.Item(4, 10).Value = 5.183
$ws.Cells.Item(4, 11).Value = 7.121
$ws.Cells.Item(4, 12).Value = 8.542999999999999
$ws.Cells.Item(4, 13).Value = 8.839
$ws.Cells.Item(4, 14).Value = 2.648
$ws.Cells.Item(4, 15).Value = 7.519
$ws.Cells.Item(4, 16).Value = 10.524
$ws.Cells.Item(4, 17).Value = 6.622
$ws.Cells.Item(4, 18).Value = 0.766
$ws.Cells.Item(4, 19).Value = 0.374
$ws.Cells.Item(4, 20).Value = 107.139
$ws.Cells.Item(4, 21).Value = 21.14
$ws.Cells.Item(4, 22).Value = 6.941
$ws.Cells.Item(4, 23).Value = 13.916
$ws.Cells.Item(4, 24).Value = 7.396
$ws.Cells.Item(4, 25).Value = 1.034
$ws.Cells.Item(4, 26).Value = 15.775
$ws.Cells.Item(4, 27).Value = 6.131
$ws.Cells.Item(4, 28).Value = 5.559
$ws.Cells.Item(4, 29).Value = 6.512
$ws.Cells.Item(4, 30).Value = 8.811
$ws.Cells.Item(4, 31).Value = 0.671
$ws.Cells.Item(4, 32).Value = 29.295
$ws.Cells.Item(4, 33).Value = 3.791
$ws.Cells.Item(4, 34).Value = 8.678000000000001
$ws.Cells.Item(5, 1).Value = 45093.52777777778
$ws.Cells.Item(5, 2).Value = 8.98
$ws.Cells.Item(5, 3).Value = 6.69
$ws.Cells.Item(5, 4).Value = 0.82
$ws.Cells.Item(5, 5).Value = 19.85
$ws.Cells.Item(5, 6).Value = 15.87
$ws.Cells.Item(5, 7).Value = 6.8
$ws.Cells.Item(5, 8).Value = 28.56
$ws.Cells.Item(5, 9).Value = 11.01
$ws.Cells.Item(5, 10).Value = 4.95
$ws.Cells.Item(5, 11).Value = 6.89
$ws.Cells.Item(5, 12).Value = 8.1
$ws.Cells.Item(5, 13).Value = 8.41
$ws.Cells.Item(5, 14).Value = 2.49
$ws.Cells.Item(5, 15).Value = 7.14
$ws.Cells.Item(5, 16).Value = 10.02
$ws.Cells.Item(5, 17).Value = 6.22
$ws.Cells.Item(5, 18).Value = 0.6
$ws.Cells.Item(5, 19).Value = 0.32
$ws.Cells.Item(5, 20).Value = 101.4
$ws.Cells.Item(5, 21).Value = 20
$ws.Cells.Item(5, 22).Value = 6.59
$ws.Cells.Item(5, 23).Value = 13.23
$ws.Cells.Item(5, 24).Value = 7.02
$ws.Cells.Item(5, 25).Value = 0.96
$ws.Cells.Item(5, 26).Value = 14.19
$ws.Cells.Item(5, 27).Value = 5.82
$ws.Cells.Item(5, 28).Value = 5.25
$ws.Cells.Item(5, 29).Value = 6.16
$ws.Cells.Item(5, 30).Value = 8.369999999999999
$ws.Cells.Item(5, 31).Value = 0.5
$ws.Cells.Item(5, 32).Value = 26.12
$ws.Cells.Item(5, 33).Value = 3.63
$ws.Cells.Item(5, 34).Value = 8.24

# 3) Widen a subset of data columns from 7 to 8 characters (custom accuracy formatting)
$ws.Columns("B:B").ColumnWidth = 7.17
$ws.Columns("C:C").ColumnWidth = 7.17
$ws.Columns("G:G").ColumnWidth = 7.17
$ws.Columns("K:K").ColumnWidth = 7.17
$ws.Columns("L:L").ColumnWidth = 7.17
$ws.Columns("M:M").ColumnWidth = 7.17
$ws.Columns("O:O").ColumnWidth = 7.17
$ws.Columns("Q:Q").ColumnWidth = 7.17
$ws.Columns("X:X").ColumnWidth = 7.17
$ws.Columns("AA:AA").ColumnWidth = 7.17
$ws.Columns("AB:AB").ColumnWidth = 7.17
$ws.Columns("AC:AC").ColumnWidth = 7.17
$ws.Columns("AD:AD").ColumnWidth = 7.17
$ws.Columns("AH:AH").ColumnWidth = 7.17
